$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 entirely (Pleasant, plot 22, Seed, Warm, 2019-12-19 to 2021-04-02 entry),
# which was a duplicate/erroneous row the author removed. This shifts all subsequent
# rows up by one.
$ws.Rows(7).Delete()

# Leave selection on the row that is now row 7 (mirrors typical post-delete state).
$ws.Range("A7:XFD7").Select()
